$wb = $excel.ActiveWorkbook

# The "O-Weapon" sheet was the previously active/selected sheet; its cached
# selection moves to I8 once the user navigates away to the new sheet, and it
# is no longer the selected tab.
$oWeapon = $wb.Worksheets.Item("O-Weapon")
$oWeapon.Range("I8").Select()

# Add the new "Status Effect" sheet at the very end of the workbook (after
# the last existing sheet, "PC All Heal").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Status Effect"

# Header row, matching the other "Single"-style battle-log sheets.
$ws.Range("A1").Value = "Index"
$ws.Range("B1").Value = "NAME"
$ws.Range("C1").Value = "ROLE"
$ws.Range("D1").Value = "LIVES"
$ws.Range("E1").Value = "Position"
$ws.Range("F1").Value = "Initiative"
$ws.Range("G1").Value = "CURRENT HP"
$ws.Range("H1").Value = "CURRENT STR"
$ws.Range("I1").Value = "CURRENT AGL"
$ws.Range("J1").Value = "CURRENT MANA"
$ws.Range("K1").Value = "CURRENT DEF"
$ws.Range("L1").Value = "COMMAND"
$ws.Range("M1").Value = "TARGET"
$ws.Range("N1").Value = "Stoned"
$ws.Range("O1").Value = "Cursed"
$ws.Range("P1").Value = "Blinded"
$ws.Range("Q1").Value = "Stunned"
$ws.Range("R1").Value = "Asleep"
$ws.Range("S1").Value = "Paralyzed"
$ws.Range("T1").Value = "Poisoned"
$ws.Range("U1").Value = "Confused"
$ws.Range("V1").Value = "ACTIONS TAKEN"

# Row 2 - a Player entry ("Zappo") that cast Cure on itself.
$ws.Range("A2").Formula = "=B2"
$ws.Range("B2").Value = "Zappo"
$ws.Range("C2").Value = "Player"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("L2").Value = "Cure"
$ws.Range("M2").Value = "Zappo"

# Row 3 - an Enemy entry ("Ghoul") targeted with the new "ParaNail" command.
$ws.Range("A3").Formula = "=B3"
$ws.Range("B3").Value = "Ghoul"
$ws.Range("C3").Value = "Enemy"
$ws.Range("D3").Value = 1
$ws.Range("L3").Value = "ParaNail"
$ws.Range("M3").Value = "Zappo"

# Restore the view: new sheet becomes the active tab/sheet with L4 selected.
$ws.Range("L4").Select()
$ws.Activate()
